$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the ParticipantsTab query (cell B2) with the new/fixed Cypher query.
# The old query text is removed from the workbook entirely (its shared string
# becomes unused and drops out), and this new query text is the one now used.
$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diagWHERE s.study_name in ["LCCC 1108: Development of a Tumor Molecular Analyses Program and Its Use to Support Treatment Decisions (UNCseqTM)"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

$ws.Range("B2").Value = $newQuery

# The last saved selection on the sheet moved from A2 to C3.
$ws.Range("C3").Select() | Out-Null
